$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting (style) of the existing "DONE" column marker (I9, which
# already holds the "X" shared string with the centered/size-20 style) onto
# the five newly completed route rows (10-14), then stamp each with "X".
$ws.Range("I9").Copy()
$ws.Range("I10:I14").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("I10").Value = "X"
$ws.Range("I11").Value = "X"
$ws.Range("I12").Value = "X"
$ws.Range("I13").Value = "X"
$ws.Range("I14").Value = "X"

# Reflect the author's last selection in the sheet view.
$ws.Range("I14").Select()
